# Table_1_panel_B.xlsx — update the j_anycomment row (row 22/23) with
# corrected N and Mean/SE figures for the Female ("(2)") column's updated
# sample, and keep the numbers stored as text (matching how every other
# numeric-looking value in this table is stored) so trailing zeros such
# as "0.200" are preserved instead of collapsing to "0.2".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    # Drop back to the sheet's default (unstyled) cell style so we don't
    # leave a stray "text format" style applied to the cell.
    $rng.Style = "Normal"
}

# Row 22 (j_anycomment): N and Mean/SE for Male (1) stay the same; the
# Female (2) N/Mean, and the (1) column's N/Mean, were corrected.
Set-TextValue "B22" "15365"
Set-TextValue "C22" "0.356"
Set-TextValue "D22" "3625"
Set-TextValue "E22" "0.200"

# Row 23: the Female SE for j_anycomment was corrected too.
Set-TextValue "E23" "[0.007]"
